# Weekly refresh of "Hortaliza, Mapocho Venta Directa de Santiago - Haba" data.
# The D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), O (Origen) and P (Precio $/Kg) values for
# rows 2-13 are rotated to the values that belonged to another row before
# the edit (a weekly shift of the underlying dataset).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row, taken from the diff (row => D,J,K,L,M,P,O)
$updates = @{
    2  = @{ D = 44446; J = 15; K = 13000; L = 13000; M = 13000; P = 520; O = "Provincia de Limarí" }
    3  = @{ D = 44453; J = 55; K = 14000; L = 15000; M = 14455; P = 578; O = "Provincia de Limarí" }
    4  = @{ D = 44467; J = 35; K = 12000; L = 12000; M = 12000; P = 480; O = "Provincia de Limarí" }
    5  = @{ D = 44340; J = 25; K = 15000; L = 15000; M = 15000; P = 600; O = "Provincia de Limarí" }
    6  = @{ D = 44432; J = 15; K = 14000; L = 14000; M = 14000; P = 560; O = "Provincia del Elquí" }
    7  = @{ D = 44425; J = 25; K = 14000; L = 14000; M = 14000; P = 560; O = "Provincia de Limarí" }
    8  = @{ D = 44418; J = 12; K = 15000; L = 15000; M = 15000; P = 600; O = "Provincia de Limarí" }
    9  = @{ D = 44449; J = 30; K = 16000; L = 16000; M = 16000; P = 640; O = "Provincia de Limarí" }
    10 = @{ D = 44376; J = 15; K = 12000; L = 12000; M = 12000; P = 480; O = "Provincia de Limarí" }
    11 = @{ D = 44421; J = 20; K = 15000; L = 15000; M = 15000; P = 600; O = "Provincia de Limarí" }
    12 = @{ D = 44435; J = 15; K = 14000; L = 14000; M = 14000; P = 560; O = "Provincia de Limarí" }
    13 = @{ D = 44435; J = 15; K = 14000; L = 14000; M = 14000; P = 560; O = "Provincia del Elquí" }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]

    $ws.Cells.Item($row, 4).Value2  = $vals.D   # D: Fecha
    $ws.Cells.Item($row, 10).Value2 = $vals.J   # J: Volumen
    $ws.Cells.Item($row, 11).Value2 = $vals.K   # K: Precio minimo
    $ws.Cells.Item($row, 12).Value2 = $vals.L   # L: Precio maximo
    $ws.Cells.Item($row, 13).Value2 = $vals.M   # M: Precio promedio ponderado
    $ws.Cells.Item($row, 15).Value2 = $vals.O   # O: Origen
    $ws.Cells.Item($row, 16).Value2 = $vals.P   # P: Precio $/Kg
}
